$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3730.6333
$ws.Range("I74").Value = 4086.2632
$ws.Range("J74").Value = 3116.3635
$ws.Range("K74").Value = 4086.2632
$ws.Range("L74").Value = 3116.3635
$ws.Range("M74").Value = -3150.2632
$ws.Range("N74").Value = -4988.363499999999
$ws.Range("H76").Value = 2665.5715
$ws.Range("I76").Value = 2670.889
$ws.Range("J76").Value = 2656
$ws.Range("K76").Value = 2670.889
$ws.Range("L76").Value = 2656
$ws.Range("M76").Value = -2355.889
$ws.Range("N76").Value = -3286
$ws.Range("H77").Value = 3730.6333
$ws.Range("I77").Value = 4086.2632
$ws.Range("J77").Value = 3116.3635
$ws.Range("K77").Value = 20431.316
$ws.Range("L77").Value = 15581.8175
$ws.Range("M77").Value = -15751.316
$ws.Range("N77").Value = -24941.8175
$ws.Range("H79").Value = 2665.5715
$ws.Range("I79").Value = 2670.889
$ws.Range("J79").Value = 2656
$ws.Range("K79").Value = 2670.889
$ws.Range("L79").Value = 2656
$ws.Range("M79").Value = -1578.889
$ws.Range("N79").Value = -4840
$ws.Range("H98").Value = 51941.45
$ws.Range("I98").Value = 84917.414
$ws.Range("J98").Value = 2477.5
$ws.Range("K98").Value = 84917.414
$ws.Range("L98").Value = 2477.5
$ws.Range("M98").Value = -83419.414
$ws.Range("N98").Value = -5473.5
$ws.Range("H112").Value = 1576.5641
$ws.Range("J112").Value = 1754
$ws.Range("L112").Value = 5262
$ws.Range("N112").Value = -7478
$ws.Range("H122").Value = 51941.45
$ws.Range("I122").Value = 84917.414
$ws.Range("J122").Value = 2477.5
$ws.Range("K122").Value = 254752.242
$ws.Range("L122").Value = 7432.5
$ws.Range("M122").Value = -252302.242
$ws.Range("N122").Value = -12332.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6485.5713
$ws.Range("I32").Value = 2752.5098
$ws.Range("J32").Value = 44562.8
$ws.Range("K32").Value = 2752.5098
$ws.Range("L32").Value = 44562.8
$ws.Range("M32").Value = -2465.5098
$ws.Range("N32").Value = -45136.8
$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1288
$ws.Range("N61").ClearContents()
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1950
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2402
$ws.Range("I99").Value = 2214.2354
$ws.Range("K99").Value = 2214.2354
$ws.Range("M99").Value = -716.2354
$ws.Range("H122").Value = 1557.6364
$ws.Range("I122").Value = 1121.8
$ws.Range("J122").Value = 1920.8334
$ws.Range("K122").Value = 3365.4
$ws.Range("L122").Value = 5762.5002
$ws.Range("M122").Value = -915.3999999999996
$ws.Range("N122").Value = -10662.5002
$ws.Range("H126").Value = 2402
$ws.Range("I126").Value = 2214.2354
$ws.Range("K126").Value = 6642.706200000001
$ws.Range("M126").Value = -4172.706200000001
$ws.Range("H132").Value = 2322.9412
$ws.Range("I132").Value = 1789.1052
$ws.Range("J132").Value = 2999.1333
$ws.Range("K132").Value = 5367.3156
$ws.Range("L132").Value = 8997.3999
$ws.Range("M132").Value = -2837.3156
$ws.Range("N132").Value = -14057.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 91.583336
$ws.Range("I12").Value = 16.666666
$ws.Range("J12").Value = 116.55556
$ws.Range("K12").Value = 49.999998
$ws.Range("L12").Value = 349.66668
$ws.Range("M12").Value = 123.000002
$ws.Range("N12").Value = -695.66668
$ws.Range("H122").Value = 2381742.5
$ws.Range("I122").Value = 682
$ws.Range("J122").Value = 4167537.8
$ws.Range("K122").Value = 6138
$ws.Range("L122").Value = 37507840.2
$ws.Range("M122").Value = -3688
$ws.Range("N122").Value = -37512740.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2945.6924
$ws.Range("I102").Value = 1888.2222
$ws.Range("K102").Value = 1888.2222
$ws.Range("M102").Value = -266.2221999999999
$ws.Range("H113").Value = 14012.5
$ws.Range("I113").Value = 34400
$ws.Range("J113").Value = 1780
$ws.Range("K113").Value = 34400
$ws.Range("L113").Value = 1780
$ws.Range("M113").Value = -32230
$ws.Range("N113").Value = -6120
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H126").Value = 2215.7778
$ws.Range("I126").Value = 1875.5
$ws.Range("J126").Value = 2488
$ws.Range("K126").Value = 5626.5
$ws.Range("L126").Value = 7464
$ws.Range("M126").Value = -3156.5
$ws.Range("N126").Value = -12404
$ws.Range("H132").Value = 9774.857
$ws.Range("I132").Value = 11713.6
$ws.Range("J132").Value = 4928
$ws.Range("K132").Value = 35140.8
$ws.Range("L132").Value = 14784
$ws.Range("M132").Value = -32610.8
$ws.Range("N132").Value = -19844

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8956.25
$ws.Range("I68").Value = 26625
$ws.Range("J68").Value = 3066.6667
$ws.Range("K68").Value = 26625
$ws.Range("L68").Value = 3066.6667
$ws.Range("M68").Value = -25876
$ws.Range("N68").Value = -4564.6667
$ws.Range("H71").Value = 8956.25
$ws.Range("I71").Value = 26625
$ws.Range("J71").Value = 3066.6667
$ws.Range("K71").Value = 133125
$ws.Range("L71").Value = 15333.3335
$ws.Range("M71").Value = -129381
$ws.Range("N71").Value = -22821.3335
$ws.Range("H122").Value = 2968.318
$ws.Range("I122").Value = 2947.2942
$ws.Range("J122").Value = 3039.8
$ws.Range("K122").Value = 8841.882599999999
$ws.Range("L122").Value = 9119.400000000001
$ws.Range("M122").Value = -6391.882599999999
$ws.Range("N122").Value = -14019.4
$ws.Range("H132").Value = 2102.054
$ws.Range("I132").Value = 1455.6207
$ws.Range("J132").Value = 4445.375
$ws.Range("K132").Value = 4366.8621
$ws.Range("L132").Value = 13336.125
$ws.Range("M132").Value = -1836.8621
$ws.Range("N132").Value = -18396.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1471.4286
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 1675
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 5025
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -9925
